# Applies the Jan 29 2024 cryptos.xlsx data refresh (GitHub Actions bot update).
# Updates Price (column D) and Volume(1h) (column E) figures for the crypto
# table, and fixes the ordering/labels of a few rows (FirstDigitalUSD/Filecoin
# and Maker/VeChain/EnergySwap/ApeXProtocol) whose coin/link values had drifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.087.82"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.258.13"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'307.47"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'96.96"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'34.69"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("D11").Value = "'0.0817"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "'6.80"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "2.608.70"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "'14.57"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "2.262.65"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "41.964.49"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "'67.40"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "'235.47"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D27").Value = "'23.44"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").Value = "'36.80"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "'164.37"
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.20"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "'3.08"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "'17.46"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").Value = "'4.13"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.939.64"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0281"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'18.75"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'2.19"
$ws.Range("E45").Value = "  -9.01%  "
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "'9.68"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "'53.98"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "2.481.00"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "'71.45"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'91.49"
$ws.Range("E51").Value = "  -0.61%  "
